$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Add the six new archival-description rows (MCH217-1 .. MCH217-6).
# Columns: A=identifier, C=title, D=date_s(blank), E=levelOfDescription,
#          F=extentAndMedium, G=notes, H=file_path(blank). Column B is left
#          untouched (no alternativeIdentifiers data for any of these rows).
# ---------------------------------------------------------------------------

$rows = @(
    @{ id = "MCH217-1"; title = "LETTES & DOCUMENTS PERTAINING TO SA'S NUCLEAR SCIENTISTS, DOCUMENTS RELATING TO SA M..... TO GERMANY"; level = "Series"; extent = "1 Box"; notes = "LOCATION: 24D | GRAP COUNT NUMER: NONE" },
    @{ id = "MCH217-2"; title = "DOCUMENTS RELATING TOO PROJEKT ADVOKAAT"; level = "Series"; extent = "1 Box"; notes = "LOCATION: 24D | GRAP COUNT NUMER: NONE" },
    @{ id = "MCH217-3"; title = "LETTERS & DOCUMENTS RELATING TOO SOUTH AFRICANS SCIENTISTS VISIT TO EUROPE, LETTERS & DOCUMENTS RELATING TO VISIT TOO SA & SOUTH WEST AFRICA BY GERMAN GENERAL"; level = "Series"; extent = "1 Box"; notes = "LOCATION: 24D | GRAP COUNT NUMER: NONE" },
    @{ id = "MCH217-4"; title = "PROKECT ADVOKAAT 2"; level = "Series"; extent = "1 Box"; notes = "LOCATION: 24D | GRAP COUNT NUMER: NONE" },
    @{ id = "MCH217-5"; title = "PROKECT ADVOKAAT 3 "; level = "Series"; extent = "1 Box"; notes = "LOCATION: 24E | GRAP COUNT NUMER: NONE" },
    @{ id = "MCH217-6"; title = "KOEBERG FUEL URANIUM, SA NUCLEAR INDUSTRY, NUCLEAR- OTHER COUNTRIES, ENRIICHMENT, NUCLEAR COLLABARATION, NUCLEAR- MILITARY, NUCLEAR AXIS"; level = "Series"; extent = "1 Box"; notes = "LOCATION: 24E | GRAP COUNT NUMER: NONE" }
)

$r = 2
foreach ($row in $rows) {
    $ws.Range("A$r").Value = $row.id
    $ws.Range("C$r").Value = $row.title
    $ws.Range("E$r").Value = $row.level
    $ws.Range("F$r").Value = $row.extent
    $ws.Range("G$r").Value = $row.notes
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# Apply the data-row styling (Calibri 10, automatic/theme text colour) to
# every populated/blank data cell in A,C,D,E,F,G,H for rows 2-7 (column B is
# skipped, matching the source data which has no alternativeIdentifiers).
# ---------------------------------------------------------------------------

$dataRange = $ws.Range("A2:A7,C2:E7,G2:H7")
foreach ($area in $dataRange.Areas) {
    $area.Font.Name = "Calibri"
    $area.Font.Size = 10
    $area.Font.ThemeColor = 1
}

$extentRange = $ws.Range("F2:F7")
$extentRange.Font.Name = "Calibri"
$extentRange.Font.Size = 10
$extentRange.Font.ThemeColor = 1
$extentRange.WrapText = $false

Write-Host "MCH217 rows added"
